# Auto-generated edit script: update cached market-price / profit values
# per the commit's scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3630.5417
$ws.Range("I113").Value = 3164.3845
$ws.Range("J113").Value = 4181.4546
$ws.Range("K113").Value = 3164.3845
$ws.Range("L113").Value = 4181.4546
$ws.Range("M113").Value = 89.61549999999988
$ws.Range("N113").Value = -10689.4546
$ws.Range("H132").Value = 7701.067
$ws.Range("I132").Value = 1806.6154
$ws.Range("J132").Value = 18647.904
$ws.Range("K132").Value = 5419.8462
$ws.Range("L132").Value = 55943.712
$ws.Range("M132").Value = -2889.8462
$ws.Range("N132").Value = -61003.712
$ws.Range("H137").Value = 10101216
$ws.Range("I137").Value = 589806.6
$ws.Range("K137").Value = 1769419.8
$ws.Range("M137").Value = -1766869.8
$ws.Range("H138").Value = 2494.8706
$ws.Range("I138").Value = 1743.3182
$ws.Range("J138").Value = 2757.3174
$ws.Range("K138").Value = 5229.9546
$ws.Range("L138").Value = 8271.9522
$ws.Range("M138").Value = -89.95460000000003
$ws.Range("N138").Value = -18551.9522
$ws.Range("H141").Value = 3450.3333
$ws.Range("I141").Value = 3406.4
$ws.Range("K141").Value = 10219.2
$ws.Range("M141").Value = -5039.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1943977.9
$ws.Range("I2").Value = 1943977.9
$ws.Range("K2").Value = 1943977.9
$ws.Range("M2").Value = -1943864.9
$ws.Range("H32").Value = 19443.611
$ws.Range("I32").Value = 20241
$ws.Range("J32").Value = 16689
$ws.Range("K32").Value = 20241
$ws.Range("L32").Value = 16689
$ws.Range("M32").Value = -19954
$ws.Range("N32").Value = -17263
$ws.Range("H43").Value = 18684.2
$ws.Range("J43").Value = 18684.2
$ws.Range("L43").Value = 18684.2
$ws.Range("N43").Value = -19310.2
$ws.Range("H63").Value = 3111
$ws.Range("I63").Value = 2325
$ws.Range("J63").Value = 3897
$ws.Range("K63").Value = 2325
$ws.Range("L63").Value = 3897
$ws.Range("M63").Value = -1639
$ws.Range("N63").Value = -5269
$ws.Range("H66").Value = 3111
$ws.Range("I66").Value = 2325
$ws.Range("J66").Value = 3897
$ws.Range("K66").Value = 11625
$ws.Range("L66").Value = 19485
$ws.Range("M66").Value = -8193
$ws.Range("N66").Value = -26349
$ws.Range("H116").Value = 1943977.9
$ws.Range("I116").Value = 1943977.9
$ws.Range("K116").Value = 1943977.9
$ws.Range("M116").Value = -1941683.9
$ws.Range("H122").Value = 3938.2
$ws.Range("I122").Value = 2147.25
$ws.Range("K122").Value = 6441.75
$ws.Range("M122").Value = -3991.75
$ws.Range("H126").Value = 5640
$ws.Range("I126").Value = 5640
$ws.Range("K126").Value = 16920
$ws.Range("M126").Value = -14450
$ws.Range("H132").Value = 12152.629
$ws.Range("I132").Value = 15709.366
$ws.Range("K132").Value = 47128.098
$ws.Range("M132").Value = -44598.098

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1943977.9
$ws.Range("I3").Value = 1943977.9
$ws.Range("K3").Value = 1943977.9
$ws.Range("M3").Value = -1943863.9
$ws.Range("H20").Value = 3094.0386
$ws.Range("J20").Value = 3329.9
$ws.Range("L20").Value = 3329.9
$ws.Range("N20").Value = -3823.9
$ws.Range("H107").Value = 1311.2273
$ws.Range("I107").Value = 1185.3334
$ws.Range("K107").Value = 1185.3334
$ws.Range("M107").Value = 734.6666
$ws.Range("H140").Value = 99997
$ws.Range("J140").Value = 99997
$ws.Range("L140").Value = 99997
$ws.Range("N140").Value = -110357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4815.74
$ws.Range("J31").Value = 8140.2
$ws.Range("L31").Value = 8140.2
$ws.Range("N31").Value = -8730.200000000001
$ws.Range("H33").Value = 10031
$ws.Range("I33").Value = 10031
$ws.Range("K33").Value = 10031
$ws.Range("M33").Value = -9652
$ws.Range("H34").Value = 4815.74
$ws.Range("J34").Value = 8140.2
$ws.Range("L34").Value = 8140.2
$ws.Range("N34").Value = -8544.200000000001
$ws.Range("H86").Value = 8275491
$ws.Range("I86").Value = 22736820
$ws.Range("J86").Value = 11874.857
$ws.Range("K86").Value = 22736820
$ws.Range("L86").Value = 11874.857
$ws.Range("M86").Value = -22735697
$ws.Range("N86").Value = -14120.857
$ws.Range("H89").Value = 8275491
$ws.Range("I89").Value = 22736820
$ws.Range("J89").Value = 11874.857
$ws.Range("K89").Value = 113684100
$ws.Range("L89").Value = 59374.285
$ws.Range("M89").Value = -113678484
$ws.Range("N89").Value = -70606.285
$ws.Range("H105").Value = 1338584
$ws.Range("I105").Value = 2274378.5
$ws.Range("K105").Value = 2274378.5
$ws.Range("M105").Value = -2272631.5
$ws.Range("H122").Value = 5093.4116
$ws.Range("I122").Value = 3058.8
$ws.Range("K122").Value = 9176.400000000001
$ws.Range("M122").Value = -6726.400000000001
$ws.Range("H132").Value = 11920568
$ws.Range("I132").Value = 13346652
$ws.Range("J132").Value = 36533
$ws.Range("K132").Value = 40039956
$ws.Range("L132").Value = 109599
$ws.Range("M132").Value = -40037426
$ws.Range("N132").Value = -114659
$ws.Range("H141").Value = 79320.09
$ws.Range("J141").Value = 89689.28999999999
$ws.Range("L141").Value = 89689.28999999999
$ws.Range("N141").Value = -100049.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 280.53845
$ws.Range("I7").Value = 249.75
$ws.Range("J7").Value = 329.8
$ws.Range("K7").Value = 749.25
$ws.Range("L7").Value = 989.4000000000001
$ws.Range("M7").Value = -637.25
$ws.Range("N7").Value = -1213.4
$ws.Range("H82").Value = 5888.8335
$ws.Range("I82").Value = 3333
$ws.Range("K82").Value = 9999
$ws.Range("M82").Value = -9593
$ws.Range("H85").Value = 5888.8335
$ws.Range("I85").Value = 3333
$ws.Range("K85").Value = 9999
$ws.Range("M85").Value = -8595
$ws.Range("H98").Value = 1272.875
$ws.Range("J98").Value = 1319
$ws.Range("L98").Value = 3957
$ws.Range("N98").Value = -6953
$ws.Range("H129").Value = 1186.8182
$ws.Range("J129").Value = 1931.6666
$ws.Range("L129").Value = 5794.9998
$ws.Range("N129").Value = -15794.9998
$ws.Range("H132").Value = 3841.5715
$ws.Range("I132").Value = 1079.8572
$ws.Range("K132").Value = 9718.7148
$ws.Range("M132").Value = -7188.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 51888.375
$ws.Range("J123").Value = 51888.375
$ws.Range("L123").Value = 51888.375
$ws.Range("N123").Value = -56788.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5830.4707
$ws.Range("I46").Value = 3741.8572
$ws.Range("K46").Value = 3741.8572
$ws.Range("M46").Value = -3553.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 785476.0600000001
$ws.Range("I100").Value = 908796.8
$ws.Range("K100").Value = 1817593.6
$ws.Range("M100").Value = -1817052.6
$ws.Range("H122").Value = 3629.7097
$ws.Range("I122").Value = 3437.36
$ws.Range("K122").Value = 10312.08
$ws.Range("M122").Value = -7862.08
$ws.Range("H136").Value = 7638.89
$ws.Range("I136").Value = 3128.1892
$ws.Range("J136").Value = 10288.032
$ws.Range("K136").Value = 9384.567599999998
$ws.Range("L136").Value = 30864.096
$ws.Range("M136").Value = -6834.567599999998
$ws.Range("N136").Value = -35964.096
